# Apply the "replace todos with actual descriptions" edit.
$wb = $excel.ActiveWorkbook

# --- Rename the third worksheet ("Include from Observation Valu 2" -> "Include from LOINC") ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Include from LOINC"

# --- Update Metadata sheet values ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2021-12-08T13:27:26-05:00"
$wsMeta.Range("B13").Value = "Frequency values for SPLASCH observations"

# --- Rebuild sheet3 as the "Include from LOINC" concept table ---
# First propagate the existing cell formatting onto the cells that will be
# newly populated, so the new rows keep the same look as the existing ones.
$ws3.Range("A1").Copy($ws3.Range("B1"))       # header style (bold) -> B1
$ws3.Range("B3").Copy($ws3.Range("B2"))       # body style -> B2
$ws3.Range("A4:B4").Copy($ws3.Range("A5:B5")) # body style -> row 5
$ws3.Range("A4:B4").Copy($ws3.Range("A6:B6")) # body style -> row 6
$ws3.Range("A4:B4").Copy($ws3.Range("A7:B7")) # body style -> row 7 (blank separator)
$ws3.Range("A4:B4").Copy($ws3.Range("A8:B8")) # body style -> row 8

# Now write the actual values.
$ws3.Range("A1").Value = "Concept"
$ws3.Range("B1").Value = "Description"

$ws3.Range("A2").Value = "LA6270-8"
$ws3.Range("B2").Value = "Never"

$ws3.Range("A3").Value = "LA10066-1"
$ws3.Range("B3").Value = "Rarely"

$ws3.Range("A4").Value = "LA10082-8"
$ws3.Range("B4").Value = "Sometimes"

$ws3.Range("A5").Value = "LA10044-8"
$ws3.Range("B5").Value = "Often"

$ws3.Range("A6").Value = "LA9933-8"
$ws3.Range("B6").Value = "Always"

$ws3.Range("A7").Value = ""
$ws3.Range("B7").Value = ""

$ws3.Range("A8").Value = "System URI"
$ws3.Range("B8").Value = "http://loinc.org"
